$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7942447662353516
$ws.Range("B1").Value = 1.491762042045593
$ws.Range("C1").Value = 5.774016380310059
$ws.Range("D1").Value = 3.13930082321167
$ws.Range("E1").Value = 1.475869536399841
